# Apply the "Updated cryptos list" data refresh to Sheet1.
# Column D ("Price") and E ("Volume(1h)") are stored as plain text in the
# workbook (prices use "." as a thousands separator, e.g. "24.780.98", and
# percentages keep their padding spaces), so any value that Excel could
# otherwise auto-parse as a number is written with the cell temporarily
# forced to Text format ("@") and then returned to General -- this keeps the
# cell content literal without leaving a lasting number-format change behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.780.98"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "1.685.06"
$ws.Range("E3").Value = "  -1.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.73%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.79"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3932"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3975"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.20%  "

$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.428"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.81%  "

$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.002"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.00"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -3.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08682"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.26"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -4.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.327"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.805"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -4.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001322"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -3.15%  "

$ws.Range("D17").Value = "1.640.35"
$ws.Range("E17").Value = "  -3.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.24"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -3.23%  "

$ws.Range("E19").Value = "  -1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.15"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -2.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.160"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.14"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.15%  "

$ws.Range("D24").Value = "24.834.03"
$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.383"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "23.87"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +2.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.784"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -7.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.70"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -3.28%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.99"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +3.47%  "

$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.777"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.609"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +18.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.865"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -8.29%  "

$ws.Range("D33").Value = "1.815.13"
$ws.Range("E33").Value = "  -4.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08477"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -4.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03089"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -2.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.019"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -4.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.967"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -4.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2818"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09577"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +3.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.52"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -3.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8000"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -6.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.74"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -3.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.458"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.72"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -5.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7199"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -3.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.594"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -4.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.193"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08750"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +5.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.344"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -3.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.54"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.59%  "
